# The "municipio-nombre" column (column E in the metadata sheet) was being
# described as a measure; it is now re-classified as a curated dimension,
# matching the pattern already used by "provincia-nombre" (F) and
# "comarca-nombre" (I).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"
